# Update gh-pages output data (合肥-漫展信息.xlsx) to match the newly
# generated scrape results.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 5185   # was 5183
$wsExhibit.Range("F7").Value = 66     # was 64
$wsExhibit.Range("F9").Value = 342    # was 339
$wsExhibit.Range("F10").Value = 6     # was 4

# --- Sheet "演出" (Shows) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("G2").Value = "不可售"   # was 180 (now no longer sellable)

# --- Sheet "全部类型" (All types, aggregated view) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("G3").Value = "不可售"    # was 180
$wsAll.Range("F9").Value = 5185        # was 5183
$wsAll.Range("F11").Value = 66         # was 64
$wsAll.Range("F14").Value = 342        # was 339
$wsAll.Range("F15").Value = 6          # was 4
